# Update cryptos list (Price and Volume(1h) columns) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.861.51"
$ws.Cells.Item(3, 4).Value = "1.813.17"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "309.13"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4977"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3884"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.09608"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.099"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "40.25"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "6.421"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.000"
$ws.Cells.Item(15, 4).Value = "1.821.13"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "7.251"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001131"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "93.32"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06597"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.15"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.936"
$ws.Cells.Item(23, 4).Value = "27.912.77"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.18"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.246"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "157.26"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.73"
$ws.Cells.Item(28, 4).Value = "2.023.52"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.407"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "128.11"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.052"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.585"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.06801"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "8.962"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02314"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2147"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "4.929"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.6235"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.000"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.144"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "13.04"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.5919"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.292"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.689"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "123.81"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.955"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.178"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06788"

$ws.Cells.Item(2, 5).Value = "  -1.15%  "
$ws.Cells.Item(3, 5).Value = "  +0.97%  "
$ws.Cells.Item(4, 5).Value = "  -0.18%  "
$ws.Cells.Item(5, 5).Value = "  -1.63%  "
$ws.Cells.Item(6, 5).Value = "  -0.09%  "
$ws.Cells.Item(7, 5).Value = "  -4.07%  "
$ws.Cells.Item(8, 5).Value = "  +1.78%  "
$ws.Cells.Item(9, 5).Value = "  +21.34%  "
$ws.Cells.Item(10, 5).Value = "  +0.07%  "
$ws.Cells.Item(11, 5).Value = "  -2.94%  "
$ws.Cells.Item(12, 5).Value = "  +2.30%  "
$ws.Cells.Item(13, 5).Value = "  -0.17%  "
$ws.Cells.Item(14, 5).Value = "  -0.48%  "
$ws.Cells.Item(15, 5).Value = "  +1.37%  "
$ws.Cells.Item(16, 5).Value = "  -0.13%  "
$ws.Cells.Item(17, 5).Value = "  +4.29%  "
$ws.Cells.Item(18, 5).Value = "  +0.24%  "
$ws.Cells.Item(19, 5).Value = "  +0.64%  "
$ws.Cells.Item(21, 5).Value = "  -0.87%  "
$ws.Cells.Item(22, 5).Value = "  -0.22%  "
$ws.Cells.Item(23, 5).Value = "  -1.13%  "
$ws.Cells.Item(24, 5).Value = "  +0.35%  "
$ws.Cells.Item(25, 5).Value = "  -0.89%  "
$ws.Cells.Item(26, 5).Value = "  -2.20%  "
$ws.Cells.Item(27, 5).Value = "  +1.37%  "
$ws.Cells.Item(28, 5).Value = "  +1.12%  "
$ws.Cells.Item(29, 5).Value = "  +3.01%  "
$ws.Cells.Item(30, 5).Value = "  +3.98%  "
$ws.Cells.Item(31, 5).Value = "  +0.41%  "
$ws.Cells.Item(32, 5).Value = "  -0.16%  "
$ws.Cells.Item(33, 5).Value = "  +0.39%  "
$ws.Cells.Item(34, 5).Value = "  -1.28%  "
$ws.Cells.Item(35, 5).Value = "  -6.90%  "
$ws.Cells.Item(36, 5).Value = "  +4.01%  "
$ws.Cells.Item(37, 5).Value = "  -0.61%  "
$ws.Cells.Item(38, 5).Value = "  +0.30%  "
$ws.Cells.Item(39, 5).Value = "  -7.25%  "
$ws.Cells.Item(40, 5).Value = "  -2.72%  "
$ws.Cells.Item(41, 5).Value = "  +1.11%  "
$ws.Cells.Item(43, 5).Value = "  -1.51%  "
$ws.Cells.Item(44, 5).Value = "  -1.47%  "
$ws.Cells.Item(45, 5).Value = "  -1.65%  "
$ws.Cells.Item(46, 5).Value = "  -5.72%  "
$ws.Cells.Item(47, 5).Value = "  -2.45%  "
$ws.Cells.Item(48, 5).Value = "  -2.72%  "
$ws.Cells.Item(49, 5).Value = "  +1.76%  "
$ws.Cells.Item(50, 5).Value = "  -4.36%  "
$ws.Cells.Item(51, 5).Value = "  +0.12%  "
